# Applies the league-table base update for rows whose id (column B) changed
# position within the sheet (Poland I Liga), matching the source XML diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowUpdates = @{}
$rowUpdates[51] = @{"B" = 5139053; "F" = "Chrobry Glogow"; "G" = "Zaglebie Sosnowiec"; "H" = 0; "I" = 0; "K" = 2.45; "L" = 3.2; "M" = 2.55; "N" = 2.7; "O" = 3.2; "P" = 2.375; "Q" = 0; "R" = 2.05; "S" = 1.75; "T" = 2.25; "U" = 1.875; "V" = 1.925; "X" = 2.2; "Z" = 0; "AA" = -0; "AB" = -1; "AC" = 0.925}
$rowUpdates[52] = @{"B" = 5139054; "F" = "GKS Tychy 71"; "G" = "Sandecja Nowy Sacz"; "H" = 2; "I" = 3; "J" = "A"; "K" = 2.15; "M" = 3.1; "N" = 2.375; "O" = 3; "P" = 3; "Q" = -0.25; "R" = 2.025; "S" = 1.775; "U" = 1.975; "V" = 1.825; "X" = -1; "Y" = 2; "Z" = -1; "AA" = 0.7749999999999999; "AB" = 0.9750000000000001; "AC" = -1}
$rowUpdates[54] = @{"B" = 5140743; "F" = "Stal Rzeszow"; "G" = "Termalica BB Nieciecza"; "I" = 2; "J" = "D"; "K" = 3; "L" = 3.3; "M" = 2.2; "N" = 2.9; "O" = 3.3; "P" = 2.25; "Q" = 0.25; "R" = 1.825; "S" = 1.975; "T" = 2.5; "U" = 1.95; "V" = 1.85; "X" = 2.3; "Y" = -1; "Z" = 0.4125; "AA" = -0.5; "AB" = 0.95}
$rowUpdates[136] = @{"B" = 5451610; "F" = "Stal Rzeszow"; "G" = "Skra Czestochowa"; "H" = 2; "J" = "H"; "K" = 1.444; "M" = 6.5; "N" = 1.333; "O" = 4.333; "P" = 8; "Q" = -1.5; "R" = 1.95; "S" = 1.85; "T" = 2.75; "U" = 1.875; "V" = 1.925; "W" = 0.333; "Y" = -1; "AA" = 0.8500000000000001; "AB" = 0.4375; "AC" = -0.5}
$rowUpdates[137] = @{"B" = 5452381; "F" = "MKS Puszcza Niepolomice"; "G" = "Chrobry Glogow"; "H" = 0; "J" = "A"; "K" = 1.571; "M" = 5; "N" = 1.4; "O" = 4.5; "P" = 6.5; "Q" = -1.25; "R" = 1.9; "S" = 1.95; "T" = 3; "U" = 2.025; "V" = 1.825; "W" = -1; "Y" = 5.5; "AA" = 0.95; "AB" = -1; "AC" = 0.825}
$rowUpdates[138] = @{"B" = 5451607; "F" = "Podbeskidzie Bielsko Biala"; "G" = "Resovia Rzeszow"; "H" = 4; "I" = 3; "K" = 1.615; "L" = 3.75; "M" = 4.75; "N" = 1.363; "P" = 6; "Q" = -1.25; "R" = 1.825; "S" = 1.975; "T" = 3.25; "U" = 1.925; "V" = 1.875; "W" = 0.363; "Z" = -0.5; "AA" = 0.4875; "AB" = 0.925; "AC" = -1}
$rowUpdates[140] = @{"B" = 5448050; "F" = "Chojniczanka Chojnice"; "G" = "GKS Katowice"; "H" = 3; "I" = 3; "J" = "D"; "K" = 2.75; "L" = 3.25; "M" = 2.375; "N" = 2.4; "O" = 3.25; "P" = 2.7; "Q" = 0; "R" = 1.8; "S" = 2.05; "T" = 2.5; "U" = 2; "V" = 1.85; "W" = -1; "X" = 2.25; "AB" = 1; "AC" = -1}
$rowUpdates[141] = @{"B" = 5448049; "F" = "LKS Lodz"; "G" = "Odra Opole"; "I" = 0; "J" = "H"; "K" = 1.571; "L" = 3.75; "M" = 5; "N" = 1.444; "O" = 4; "P" = 6; "Q" = -1; "R" = 1.775; "S" = 2.025; "T" = 2.75; "U" = 1.9; "V" = 1.9; "W" = 0.444; "X" = -1; "Z" = 0; "AA" = -0; "AB" = -1; "AC" = 0.8999999999999999}
$rowUpdates[142] = @{"B" = 5448048; "F" = "Zaglebie Sosnowiec"; "G" = "Sandecja Nowy Sacz"; "H" = 1; "I" = 1; "J" = "D"; "K" = 2.1; "L" = 3.2; "M" = 3.3; "N" = 2.1; "O" = 3.2; "P" = 3.1; "Q" = -0.25; "R" = 1.875; "S" = 1.925; "T" = 2.25; "U" = 1.85; "V" = 1.95; "X" = 2.2; "Y" = -1; "Z" = -0.5; "AA" = 0.4625; "AB" = -0.5; "AC" = 0.475}
$rowUpdates[143] = @{"B" = 5451609; "F" = "Ruch Chorzow"; "G" = "GKS Tychy 71"; "H" = 1; "I" = 0; "K" = 1.3; "L" = 5; "M" = 7.5; "N" = 1.333; "P" = 8; "Q" = -1.5; "R" = 2; "S" = 1.8; "T" = 2.5; "U" = 1.825; "V" = 1.975; "W" = 0.333; "Z" = -1; "AA" = 0.8; "AB" = -1; "AC" = 0.9750000000000001}
$rowUpdates[144] = @{"B" = 5447925; "F" = "Gornik Leczna"; "G" = "Wisla Krakow"; "H" = 0; "J" = "A"; "K" = 5.5; "L" = 4; "M" = 1.5; "N" = 4.5; "O" = 4; "P" = 1.615; "Q" = 0.75; "R" = 2.05; "S" = 1.8; "T" = 3; "X" = -1; "Y" = 0.615; "Z" = -1; "AA" = 0.8; "AB" = 0; "AC" = -0}
$rowUpdates[209] = @{"B" = 6805719; "F" = "Motor Lublin"; "G" = "Stal Rzeszow"; "H" = 3; "I" = 2; "K" = 2.3; "M" = 2.8; "N" = 2.05; "O" = 3.4; "P" = 3.2; "Q" = -0.25; "R" = 1.85; "S" = 2; "T" = 2.5; "U" = 1.875; "V" = 1.975; "W" = 1.05; "Z" = 0.8500000000000001; "AB" = 0.875; "AC" = -1}
$rowUpdates[210] = @{"B" = 6803727; "F" = "Chrobry Glogow"; "G" = "GKS Tychy 71"; "H" = 2; "I" = 1; "K" = 3.1; "M" = 2.15; "N" = 3.75; "O" = 3.6; "P" = 1.85; "Q" = 0.5; "R" = 1.925; "S" = 1.875; "T" = 2.75; "U" = 1.925; "V" = 1.875; "W" = 2.75; "Z" = 0.925; "AB" = 0.4625; "AC" = -0.5}
$rowUpdates[302] = @{"B" = 6803794; "F" = "Wisla Krakow"; "G" = "Gornik Leczna"; "H" = 4; "I" = 0; "J" = "H"; "K" = 1.4; "L" = 4.75; "M" = 7; "N" = 1.363; "O" = 4.75; "P" = 7.5; "Q" = -1.25; "R" = 1.8; "S" = 2; "T" = 2.75; "U" = 1.775; "V" = 2.025; "W" = 0.363; "X" = -1; "Z" = 0.8; "AA" = -1; "AB" = 0.7749999999999999; "AC" = -1}
$rowUpdates[303] = @{"B" = 6803793; "F" = "Odra Opole"; "G" = "Stal Rzeszow"; "H" = 1; "I" = 1; "J" = "D"; "K" = 2.05; "L" = 3.4; "M" = 3.5; "N" = 2.1; "O" = 3.4; "P" = 3.4; "Q" = -0.25; "R" = 1.825; "S" = 1.975; "T" = 2.5; "U" = 1.95; "V" = 1.85; "W" = -1; "X" = 2.4; "Z" = -0.5; "AA" = 0.4875; "AB" = -1; "AC" = 0.8500000000000001}

foreach ($r in $rowUpdates.Keys) {
    $colVals = $rowUpdates[$r]
    foreach ($col in $colVals.Keys) {
        $ws.Range("$col$r").Value2 = $colVals[$col]
    }
}

Write-Host "Updated rows:" ($rowUpdates.Keys -join ", ")